$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.673.20"
$ws.Cells.Item(2, 5).Value = "  -0.83%  "
$ws.Cells.Item(3, 4).Value = "1.585.41"
$ws.Cells.Item(3, 5).Value = "  -3.10%  "
$ws.Cells.Item(4, 5).Value = "  +0.22%  "
$ws.Cells.Item(5, 4).Value = "'206.58"
$ws.Cells.Item(5, 5).Value = "  -2.45%  "
$ws.Cells.Item(6, 5).Value = "  -3.28%  "
$ws.Cells.Item(8, 4).Value = "'22.28"
$ws.Cells.Item(8, 5).Value = "  -4.69%  "
$ws.Cells.Item(9, 5).Value = "  -1.24%  "
$ws.Cells.Item(10, 4).Value = "'0.0592"
$ws.Cells.Item(10, 5).Value = "  -3.17%  "
$ws.Cells.Item(11, 4).Value = "'0.0868"
$ws.Cells.Item(11, 5).Value = "  -1.82%  "
$ws.Cells.Item(12, 4).Value = "1.810.64"
$ws.Cells.Item(12, 5).Value = "  -3.10%  "
$ws.Cells.Item(13, 4).Value = "1.597.27"
$ws.Cells.Item(13, 5).Value = "  -2.37%  "
$ws.Cells.Item(14, 5).Value = "  -4.03%  "
$ws.Cells.Item(15, 4).Value = "'0.531"
$ws.Cells.Item(15, 5).Value = "  -5.67%  "
$ws.Cells.Item(16, 4).Value = "27.639.71"
$ws.Cells.Item(16, 5).Value = "  -1.01%  "
$ws.Cells.Item(17, 4).Value = "'63.13"
$ws.Cells.Item(17, 5).Value = "  -3.50%  "
$ws.Cells.Item(18, 4).Value = "'218.74"
$ws.Cells.Item(18, 5).Value = "  -4.45%  "
$ws.Cells.Item(19, 4).Value = "0.0₃0694"
$ws.Cells.Item(19, 5).Value = "  -3.58%  "
$ws.Cells.Item(20, 4).Value = "'7.31"
$ws.Cells.Item(20, 5).Value = "  -5.03%  "
$ws.Cells.Item(21, 5).Value = "  +0.27%  "
$ws.Cells.Item(22, 5).Value = "  -5.04%  "
$ws.Cells.Item(23, 4).Value = "'9.55"
$ws.Cells.Item(24, 4).Value = "'1.98"
$ws.Cells.Item(24, 5).Value = "  -5.04%  "
$ws.Cells.Item(25, 4).Value = "'153.61"
$ws.Cells.Item(25, 5).Value = "  -1.10%  "
$ws.Cells.Item(26, 5).Value = "  +0.29%  "
$ws.Cells.Item(27, 4).Value = "'6.73"
$ws.Cells.Item(27, 5).Value = "  -2.52%  "
$ws.Cells.Item(28, 5).Value = "  -2.98%  "
$ws.Cells.Item(29, 5).Value = "  -4.01%  "
$ws.Cells.Item(30, 5).Value = "  -2.87%  "
$ws.Cells.Item(31, 4).Value = "'0.0466"
$ws.Cells.Item(31, 5).Value = "  -3.33%  "
$ws.Cells.Item(32, 5).Value = "  -5.52%  "
$ws.Cells.Item(33, 4).Value = "1.382.46"
$ws.Cells.Item(33, 5).Value = "  -1.04%  "
$ws.Cells.Item(34, 5).Value = "  -4.94%  "
$ws.Cells.Item(35, 5).Value = "  -5.68%  "
$ws.Cells.Item(36, 4).Value = "'0.966"
$ws.Cells.Item(36, 5).Value = "  -4.60%  "
$ws.Cells.Item(37, 5).Value = "  -1.18%  "
$ws.Cells.Item(38, 4).Value = "'0.0166"
$ws.Cells.Item(38, 5).Value = "  -2.92%  "
$ws.Cells.Item(39, 5).Value = "  -3.44%  "
$ws.Cells.Item(40, 4).Value = "'0.820"
$ws.Cells.Item(40, 5).Value = "  -3.73%  "
$ws.Cells.Item(41, 5).Value = "  +0.29%  "
$ws.Cells.Item(42, 5).Value = "  -2.46%  "
$ws.Cells.Item(43, 5).Value = "  -3.31%  "
$ws.Cells.Item(44, 4).Value = "'63.84"
$ws.Cells.Item(44, 5).Value = "  -3.41%  "
$ws.Cells.Item(45, 5).Value = "  +1.50%  "
$ws.Cells.Item(46, 5).Value = "  -4.21%  "
$ws.Cells.Item(47, 4).Value = "1.721.31"
$ws.Cells.Item(47, 5).Value = "  -3.11%  "
$ws.Cells.Item(48, 4).Value = "'87.85"
$ws.Cells.Item(48, 5).Value = "  -1.05%  "
$ws.Cells.Item(49, 4).Value = "0.0₆0102"
$ws.Cells.Item(49, 5).Value = "  -1.11%  "
$ws.Cells.Item(50, 4).Value = "'0.0976"
$ws.Cells.Item(50, 5).Value = "  -5.12%  "
